# The deck's theme (ppt/theme/theme1.xml, the "Integral" design used by the
# slide master / all slides) is switched over to the stock PowerPoint
# "Office Theme" colour palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint exposes the theme's 12-colour scheme through
# Slide.ThemeColorScheme (it is shared by every slide / the slide master,
# since they all point at the same underlying theme part), so we rewrite
# each of the 12 entries to the "Office Theme" RGB values.

function HexToComRgb($hex) {
    # OOXML srgbClr is RRGGBB; COM ColorFormat.RGB / RGBColor.RGB values are
    # packed little-endian as 0x00BBGGRR, so convert accordingly.
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$scheme = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = HexToComRgb($officeThemeColors[$i - 1])
}
